# Adds the new survey response row (Harrison Driver, 2025-12-02 13:04:05)
# submitted from the Streamlit SmartScore app as row 21 of Sheet1.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("A21").Value = 'Harrison Driver_20251202_130404'
$ws.Range("B21").Value = "'"
$ws.Range("B21").Style = "Normal"
$ws.Range("C21").Value = 'Harrison Driver'
$ws.Range("D21").Value = 18
$ws.Range("E21").Value = 'Male'
$ws.Range("F21").Value = '2025-12-02 13:04:05'
$ws.Range("G21").Value = @"
{
  "portion": 1.0,
  "diet": 1.0,
  "salt": 0.2,
  "fat": 0.2,
  "natural": 0.8,
  "convenience": 1.0,
  "price": 0.8
}
"@
# Note: the SmartScore columns (I/L/O/R/U/X/AA/AD/AG) came through the
# Streamlit export as text, e.g. "0.646" rather than the numeric 0.646 used
# in earlier rows - a leading quote-prefix forces Excel to store them as
# text instead of auto-coercing to numbers (the prefix/style is cleared below).
$ws.Range("H21").Value = 'Nongshim Neoguri Spicy Seafood'
$ws.Range("I21").Value = '''0.646'
$ws.Range("J21").Value = 'Sabor a marisco, umami, picante equilibrado, buena textura, algo salado'
$ws.Range("K21").Value = 'Nissin Chow Mein Teriyaki Beef'
$ws.Range("L21").Value = '''0.543'
$ws.Range("M21").Value = 'Fácil de preparar, porción generosa, salsa suave, necesita mejoras, alto en grasa'
$ws.Range("N21").Value = 'Nongshim Shin Ramyun'
$ws.Range("O21").Value = '''0.535'
$ws.Range("P21").Value = 'Sabor intenso, picante, umami, fideos gruesos, muy alto en sodio'
$ws.Range("Q21").Value = 'Amy’s Macaroni & Cheese (frozen)'
$ws.Range("R21").Value = '''0.662'
$ws.Range("S21").Value = 'Queso real, textura casera, sin conservadores, alto en grasa, algo caro'
$ws.Range("T21").Value = 'Kraft Macaroni & Cheese Dinner'
$ws.Range("U21").Value = '''0.507'
$ws.Range("V21").Value = 'Sabor nostálgico, clásico americano, fácil, no muy nutritivo, barato'
$ws.Range("W21").Value = 'Annie’s Shells & White Cheddar'
$ws.Range("X21").Value = '''0.456'
$ws.Range("Y21").Value = 'Queso blanco real, sin colorantes, sabor casero, menos salado, buena para niños'
$ws.Range("Z21").Value = 'Wild Planet Wild Tuna Pasta Salad'
$ws.Range("AA21").Value = '''0.720'
$ws.Range("AB21").Value = 'Sabor fresco, buena proteína, saludable, porción algo pequeña'
$ws.Range("AC21").Value = 'StarKist Chicken Creations (Chicken Salad)'
$ws.Range("AD21").Value = '''0.498'
$ws.Range("AE21").Value = 'Portátil, saludable, fácil, buena textura, sabor suave'
$ws.Range("AF21").Value = 'Kitchens of India Variety Pack'
$ws.Range("AG21").Value = '''0.472'
$ws.Range("AH21").Value = 'Sabor auténtico, variedad, vegetariano, necesita arroz o pan, buena calidad'

# Clear quote-prefix styling on numeric-looking text cells so no extra style index is left behind
$ws.Range("I21").Style = "Normal"
$ws.Range("L21").Style = "Normal"
$ws.Range("O21").Style = "Normal"
$ws.Range("R21").Style = "Normal"
$ws.Range("U21").Style = "Normal"
$ws.Range("X21").Style = "Normal"
$ws.Range("AA21").Style = "Normal"
$ws.Range("AD21").Style = "Normal"
$ws.Range("AG21").Style = "Normal"
$ws.Range("B21").Style = "Normal"

# Reset row height: the multi-line JSON in G21 triggers Excel's auto row-height,
# which would otherwise stamp an explicit ht/customHeight on the row.
$ws.Rows.Item(21).AutoFit()
